# Shade the "personne" row-group of the routes table green (fill 00B050),
# matching the shading already applied to the "equipe" row-group above it.

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# wdColor value for RGB hex 00B050: R | (G << 8) | (B << 16)
# R=0x00, G=0xB0, B=0x50  =>  0x00 | (0xB0 << 8) | (0x50 << 16) = 0x50B000
$green = [int]0x50B000
$wdColorAutomatic = -16777216

# Locate the row whose first column holds the category label "personne".
$startRow = 0
for ($i = 1; $i -le $table.Rows.Count; $i++) {
    $labelCell = $table.Rows.Item($i).Cells.Item(1)
    $cellText = ($labelCell.Range.Text -replace "[\x07\r\x0c]", "").Trim()
    if ($cellText -eq "personne") {
        $startRow = $i
        break
    }
}

if ($startRow -eq 0) {
    throw "Could not locate the 'personne' row group"
}

# The category occupies 5 data rows (list all / get / add / update / delete).
$endRow = $startRow + 4

for ($r = $startRow; $r -le $endRow; $r++) {
    # Column 1 is a vertically merged cell; only reachable reliably through
    # the row's own Cells collection (Table.Cell(r, 1) misbehaves on the
    # vMerge-continuation rows).
    $cells = New-Object System.Collections.ArrayList
    [void]$cells.Add($table.Rows.Item($r).Cells.Item(1))
    [void]$cells.Add($table.Cell($r, 2))
    [void]$cells.Add($table.Cell($r, 3))

    foreach ($cell in $cells) {
        $cell.Shading.Texture = 0
        $cell.Shading.ForegroundPatternColor = $wdColorAutomatic
        $cell.Shading.BackgroundPatternColor = $green
    }
}
